# fix: product import template - add missing "Highlight" and
# "Chuong trinh dac biet" columns, inserted right after "Khuyen mai"
# (and before "Thuong hieu"), shifting the trailing detail columns
# (Thuong hieu .. Loai Len) two columns to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at O:P (old O:T -> new Q:V).
$ws.Columns("O:P").Insert()

# Populate the headers for the two newly inserted columns. The cells
# inherit the bold header style from the insert, matching the rest of
# row 1.
$ws.Cells.Item(1, 15).Value2 = "Highlight"
$ws.Cells.Item(1, 16).Value2 = "Chương trình đặc biệt"

# The column insert also stamps blank placeholder cells into row 2 for
# the two new columns (no data there, same as original row 2 which only
# goes through column N) - clear them out completely so no stray empty
# <c> nodes are written for O2/P2.
$ws.Range("O2:P2").Clear()

# Match the author's resulting selection/top-left view state.
$ws.Range("P5").Select()
